# Minutes for 2024-11-01 working group meeting
# - add new attendees Katie Harding (Freenome) and Youn Kyeong Chang (FDA)
# - open external links in a new tab (trailing-space normalization on the
#   Sanofi hyperlink display text carried over from the source edit)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (styles) from the last existing data row (89) down to the
# two new rows so that column C keeps the blue "link" look used throughout
# the sheet.
$ws.Range("A89:C89").Copy($ws.Range("A90:C90"))
$ws.Range("A89:C89").Copy($ws.Range("A91:C91"))

# Row 90: Katie Harding / Freenome
$ws.Range("A90").Value = "Katie Harding"
$ws.Range("B90").Value = "Freenome"
$ws.Range("C90").Value = "https://www.freenome.com/ "

# Row 91: Youn Kyeong Chang / FDA
$ws.Range("A91").Value = "Youn Kyeong Chang"
$ws.Range("B91").Value = "FDA"
$ws.Range("C91").Value = "https://www.fda.gov/"

# Sanofi's affiliation URL cell picks up a trailing space
$ws.Range("C87").Value = "https://www.sanofi.com/ "

# New hyperlinks for the two new affiliation URLs
$ws.Hyperlinks.Add($ws.Range("C90"), "https://www.freenome.com/", [Type]::Missing, [Type]::Missing, "https://www.freenome.com/")
$ws.Hyperlinks.Add($ws.Range("C91"), "https://www.fda.gov/", [Type]::Missing, [Type]::Missing, "https://www.fda.gov/")

# Adding the hyperlinks re-styles the target cells with Excel's built-in
# "Hyperlink" style; restore the sheet's own link style (copied from C89)
# so column C stays visually consistent.
$ws.Range("C89").Copy($ws.Range("C90"))
$ws.Range("C90").Value = "https://www.freenome.com/ "
$ws.Range("C89").Copy($ws.Range("C91"))
$ws.Range("C91").Value = "https://www.fda.gov/"

# Update the active selection to match the author's final cursor position
$ws.Range("C91").Select()

Write-Output "applied wg_attendees edits"
